$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Regresi Tanggal" update: the PERIODE_HARIAN date entered into O2 moves
# from 05/05/2023 to 07/09/2024. P2 re-derives YYYYMMDD from it, and S2 is
# re-pointed to reference P2 via a formula instead of holding a separate
# static value.
$ws.Range("O2").Value = "07/09/2024"
$ws.Range("S2").Formula = "=P2"

# Leave the selection where the author left it after making the edit.
$ws.Range("S3").Select()
